$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.929.25"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "3.116.44"
$ws.Range("E3").Value = "  +0.70%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'580.11"
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("D6").Value = "'173.34"
$ws.Range("E6").Value = "  +1.79%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  -0.30%  "
$ws.Range("D9").Value = "'6.42"
$ws.Range("E9").Value = "  -3.01%  "
$ws.Range("E10").Value = "  -1.15%  "
$ws.Range("D11").Value = "'0.480"
$ws.Range("E11").Value = "  +0.08%  "
$ws.Range("D12").Value = "'0.0000247"
$ws.Range("E12").Value = "  -0.80%  "
$ws.Range("D13").Value = "'37.34"
$ws.Range("E13").Value = "  +1.90%  "
$ws.Range("D14").Value = "'0.122"
$ws.Range("E14").Value = "  -1.48%  "
$ws.Range("D15").Value = "3.631.37"
$ws.Range("E15").Value = "  +0.67%  "
$ws.Range("D16").Value = "66.859.87"
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("D17").Value = "'7.14"
$ws.Range("E17").Value = "  -0.56%  "
$ws.Range("D18").Value = "3.115.92"
$ws.Range("E18").Value = "  +0.67%  "
$ws.Range("D19").Value = "'16.38"
$ws.Range("E19").Value = "  +1.03%  "
$ws.Range("D20").Value = "'477.80"
$ws.Range("E20").Value = "  +2.59%  "
$ws.Range("D21").Value = "'0.708"
$ws.Range("E21").Value = "  -0.92%  "
$ws.Range("D22").Value = "'7.82"
$ws.Range("E22").Value = "  +4.64%  "
$ws.Range("D23").Value = "'84.02"
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("D24").Value = "'13.19"
$ws.Range("E24").Value = "  +0.32%  "
$ws.Range("D25").Value = "'2.29"
$ws.Range("E25").Value = "  -1.81%  "
$ws.Range("D26").Value = "'10.38"
$ws.Range("E26").Value = "  +3.09%  "
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("D28").Value = "'7.95"
$ws.Range("E28").Value = "  -0.52%  "
$ws.Range("E29").Value = "  -1.36%  "
$ws.Range("E30").Value = "  +0.40%  "
$ws.Range("D31").Value = "'28.54"
$ws.Range("E31").Value = "  +1.10%  "
$ws.Range("E32").Value = "  +0.18%  "
$ws.Range("D33").Value = "0.0₃0949"
$ws.Range("E33").Value = "  -8.48%  "
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("D35").Value = "'5.84"
$ws.Range("E35").Value = "  -0.91%  "
$ws.Range("D36").Value = "'0.973"
$ws.Range("E36").Value = "  -2.80%  "
$ws.Range("D37").Value = "'47.05"
$ws.Range("E37").Value = "  -1.22%  "
$ws.Range("D38").Value = "'50.11"
$ws.Range("E38").Value = "  -0.23%  "
$ws.Range("D39").Value = "'2.05"
$ws.Range("E39").Value = "  -2.80%  "
$ws.Range("D40").Value = "'0.314"
$ws.Range("E40").Value = "  -0.92%  "
$ws.Range("E41").Value = "  +1.33%  "
$ws.Range("D42").Value = "'8.55"
$ws.Range("E42").Value = "  -1.23%  "
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").Value = "'385.52"
$ws.Range("E43").Value = "  +0.97%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "2.815.89"
$ws.Range("E44").Value = "  +0.94%  "
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").Value = "'2.58"
$ws.Range("E45").Value = "  -7.83%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "'0.0353"
$ws.Range("E46").Value = "  -1.76%  "
$ws.Range("D47").Value = "'135.93"
$ws.Range("E47").Value = "  +1.01%  "
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("D49").Value = "'25.09"
$ws.Range("E49").Value = "  +1.63%  "
$ws.Range("E50").Value = "  -1.04%  "
$ws.Range("E51").Value = "  -0.65%  "
